$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, G (rows 2-14). F column is unchanged.
$data = @{
    2  = @{ B = 0.6545652718822623;  C = 1.626987699542094;  D = 0.1496068669990043;   E = 13.86384647080068;  G = 16.29500630922404  }
    3  = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 0.7210945179870265;   E = 13.86384647080068;  G = 19.48425592650926  }
    4  = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 0.7210945179870265;   E = 0.5333859586016987; G = 6.15379541431027   }
    5  = @{ B = 0.6545652718822623;  C = 1.626987699542094;  D = 3.223369029078222;    E = 0.5333859586016987; G = 6.038307959104277  }
    6  = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 3.223369029078222;    E = 0.5333859586016987; G = 8.656069925401464  }
    7  = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 3.223369029078222;    E = 0.5333859586016987; G = 8.656069925401464  }
    8  = @{ B = 0.6545652718822623;  C = 1.626987699542094;  D = 0.7210945179870265;   E = 0.5333859586016987; G = 3.536033448013082  }
    9  = @{ B = 1.445647641019636;   C = 1.626987699542094;  D = 0.7210945179870265;   E = 0.5333859586016987; G = 4.327115817150455  }
    10 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 3.223369029078222;    E = 0.5333859586016987; G = 8.656069925401464  }
    11 = @{ B = 1.445647641019636;   C = 1.626987699542094;  D = 189.6080260415259;    E = 13.86384647080068;  G = 206.5445078528883  }
    12 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 18.71679738969934;    E = 0.5333859586016987; G = 24.14949828602258  }
    13 = @{ B = 0.6545652718822623;  C = 0.3048912486333797; D = 0.7210945179870265;   E = 0.5333859586016987; G = 2.213936997104367  }
    14 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 3.223369029078222;    E = 0.5333859586016987; G = 8.656069925401464  }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
